$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2:E51").NumberFormat = "@"

$ws.Range("D2").Value = "39.841.14"
$ws.Range("E2").Value = "  -3.18%  "
$ws.Range("D3").Value = "2.334.59"
$ws.Range("E3").Value = "  -3.94%  "
$ws.Range("E4").Value = "  -0.03%  "
$ws.Range("D5").Value = "308.92"
$ws.Range("E5").Value = "  -2.92%  "
$ws.Range("D6").Value = "83.87"
$ws.Range("E6").Value = "  -6.45%  "
$ws.Range("E7").Value = "  -2.52%  "
$ws.Range("E8").Value = "  -0.02%  "
$ws.Range("E9").Value = "  -4.17%  "
$ws.Range("D10").Value = "0.0801"
$ws.Range("E10").Value = "  -4.47%  "
$ws.Range("D11").Value = "29.76"
$ws.Range("E11").Value = "  -7.38%  "
$ws.Range("E12").Value = "  +0.36%  "
$ws.Range("D13").Value = "2.697.59"
$ws.Range("E13").Value = "  -3.79%  "
$ws.Range("D14").Value = "6.37"
$ws.Range("E14").Value = "  -5.58%  "
$ws.Range("D15").Value = "14.67"
$ws.Range("E15").Value = "  -6.37%  "
$ws.Range("D16").Value = "2.350.57"
$ws.Range("E16").Value = "  -2.93%  "
$ws.Range("D17").Value = "0.752"
$ws.Range("E17").Value = "  -3.42%  "
$ws.Range("D18").Value = "39.788.54"
$ws.Range("E18").Value = "  -3.15%  "
$ws.Range("E19").Value = "  -3.56%  "
$ws.Range("E20").Value = "  -4.26%  "
$ws.Range("D21").Value = "67.91"
$ws.Range("E21").Value = "  -6.22%  "
$ws.Range("D22").Value = "10.47"
$ws.Range("E22").Value = "  -5.50%  "
$ws.Range("D23").Value = "234.50"
$ws.Range("E23").Value = "  -0.42%  "
$ws.Range("D24").Value = "2.52"
$ws.Range("E24").Value = "  -6.30%  "
$ws.Range("E25").Value = "  -0.05%  "
$ws.Range("E26").Value = "  -3.74%  "
$ws.Range("D27").Value = "23.34"
$ws.Range("E27").Value = "  -3.60%  "
$ws.Range("D28").Value = "2.20"
$ws.Range("E28").Value = "  -1.54%  "
$ws.Range("D29").Value = "9.20"
$ws.Range("E29").Value = "  -4.59%  "
$ws.Range("D30").Value = "34.08"
$ws.Range("E30").Value = "  -1.62%  "
$ws.Range("D31").Value = "152.20"
$ws.Range("E31").Value = "  -4.21%  "
$ws.Range("E32").Value = "  +0.04%  "
$ws.Range("E33").Value = "  -4.11%  "
$ws.Range("D34").Value = "2.47"
$ws.Range("E34").Value = "  -0.87%  "
$ws.Range("D35").Value = "0.0712"
$ws.Range("E35").Value = "  -4.83%  "
$ws.Range("E36").Value = "  -1.08%  "
$ws.Range("E37").Value = "  -7.19%  "
$ws.Range("D38").Value = "0.0980"
$ws.Range("E38").Value = "  -3.02%  "
$ws.Range("D39").Value = "15.46"
$ws.Range("E39").Value = "  -9.62%  "
$ws.Range("E40").Value = "  -5.54%  "
$ws.Range("E41").Value = "  -3.72%  "
$ws.Range("D42").Value = "1.965.76"
$ws.Range("E42").Value = "  -1.55%  "
$ws.Range("D43").Value = "2.25"
$ws.Range("E43").Value = "  -3.18%  "
$ws.Range("E44").Value = "  -4.99%  "
$ws.Range("D45").Value = "17.36"
$ws.Range("E45").Value = "  -6.57%  "
$ws.Range("D46").Value = "9.41"
$ws.Range("E46").Value = "  -1.71%  "
$ws.Range("E47").Value = "  -8.49%  "
$ws.Range("D48").Value = "2.560.74"
$ws.Range("E48").Value = "  -3.91%  "
$ws.Range("D49").Value = "92.21"
$ws.Range("E49").Value = "  -2.95%  "
$ws.Range("E50").Value = "  -5.30%  "
$ws.Range("D51").Value = "49.53"
$ws.Range("E51").Value = "  -4.57%  "

$ws.Range("D2:E51").ClearFormats()
